$wb = $excel.ActiveWorkbook

# --- Update Version and Date values on the Metadata sheet ---
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B3").Value = "1.1.0"
$metadata.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# --- Re-apply the wrap/vertical-top alignment so it is flagged as "applied" ---
# (mirrors the applyAlignment="true" attribute added to the two shared cell
#  styles used across all worksheets in this workbook)
$metadata.Range("A1:B1").VerticalAlignment = -4160
$metadata.Range("A1:B1").WrapText = $true
$metadata.Range("A2:B14").VerticalAlignment = -4160
$metadata.Range("A2:B14").WrapText = $true

$includeSheetNames = @(
    "Include from FSIII",
    "Include from FSIII 2",
    "Include from FSIII 3",
    "Include from FSIII 4"
)

foreach ($name in $includeSheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A1:C1").VerticalAlignment = -4160
    $ws.Range("A1:C1").WrapText = $true

    $ws.Range("A2:C2").VerticalAlignment = -4160
    $ws.Range("A2:C2").WrapText = $true

    $ws.Range("A3:B4").VerticalAlignment = -4160
    $ws.Range("A3:B4").WrapText = $true
}
